$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column in H1 - reuse the existing header formatting (from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data values for the two data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
